$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 556, shifting existing rows 556-611 down to 558-613.
$ws.Rows("556:557").Insert()

# Row 556: new Cilantro record, same qualitative columns as the rest of the block.
$ws.Cells.Item(556, 1).Value = 10
$ws.Cells.Item(556, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(556, 3).Value = "La Araucanía"
$ws.Cells.Item(556, 4).Value = 45132
$ws.Cells.Item(556, 5).Value = 9
$ws.Cells.Item(556, 6).Value = 100112040
$ws.Cells.Item(556, 7).Value = "Cilantro"
$ws.Cells.Item(556, 8).Value = "Sin especificar"
$ws.Cells.Item(556, 9).Value = "Primera"
$ws.Cells.Item(556, 10).Value = 65
$ws.Cells.Item(556, 11).Value = 4000
$ws.Cells.Item(556, 12).Value = 4000
$ws.Cells.Item(556, 13).Value = 4000
$ws.Cells.Item(556, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(556, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(556, 16).Value = 2000
$ws.Cells.Item(556, 17).Value = 2
$ws.Cells.Item(556, 18).Value = "Hortaliza"

# Row 557: new Cilantro record.
$ws.Cells.Item(557, 1).Value = 10
$ws.Cells.Item(557, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(557, 3).Value = "La Araucanía"
$ws.Cells.Item(557, 4).Value = 45132
$ws.Cells.Item(557, 5).Value = 9
$ws.Cells.Item(557, 6).Value = 100112040
$ws.Cells.Item(557, 7).Value = "Cilantro"
$ws.Cells.Item(557, 8).Value = "Sin especificar"
$ws.Cells.Item(557, 9).Value = "Primera"
$ws.Cells.Item(557, 10).Value = 40
$ws.Cells.Item(557, 11).Value = 4000
$ws.Cells.Item(557, 12).Value = 4000
$ws.Cells.Item(557, 13).Value = 4000
$ws.Cells.Item(557, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(557, 15).Value = "Región Metropolitana"
$ws.Cells.Item(557, 16).Value = 2000
$ws.Cells.Item(557, 17).Value = 2
$ws.Cells.Item(557, 18).Value = "Hortaliza"
